$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "202.173.124.126"
$ws.Range("B11").Value = 28.3621581
$ws.Range("C11").Value = 77.2828472
$ws.Range("D11").Value = 100
$ws.Range("E11").Value = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
$ws.Range("F11").Value = "Linux armv81"
$ws.Range("G11").Value = "2025-06-25T16:41:57.591Z"
